$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 75 (pushes existing rows 75-89 down to 76-90),
# copying formatting from row 74 above (matches Excel's default Insert behaviour).
$ws.Rows.Item(75).Insert()

# The new row 75 becomes a copy of the (original) row 74 data, i.e. the data
# that used to be in row 74 before its date was updated to a newer sample date.
$ws.Cells.Item(75, 1).Value  = 1
$ws.Cells.Item(75, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(75, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(75, 4).Value  = 44673
$ws.Cells.Item(75, 5).Value  = 15
$ws.Cells.Item(75, 6).Value  = "Fruta"
$ws.Cells.Item(75, 7).Value  = 100106
$ws.Cells.Item(75, 8).Value  = "Oleaginosos"
$ws.Cells.Item(75, 9).Value  = 100106002
$ws.Cells.Item(75, 10).Value = "Palta"
$ws.Cells.Item(75, 11).Value = "Hass"
$ws.Cells.Item(75, 12).Value = "Primera"
$ws.Cells.Item(75, 13).Value = 400
$ws.Cells.Item(75, 14).Value = 18000
$ws.Cells.Item(75, 15).Value = 19000
$ws.Cells.Item(75, 16).Value = 18500
$ws.Cells.Item(75, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(75, 18).Value = "Perú"
$ws.Cells.Item(75, 19).Value = 1850
$ws.Cells.Item(75, 20).Value = 10

# Row 74 keeps its original data except for a new, more recent sample date.
$ws.Cells.Item(74, 4).Value = 44694
